# Auto-generated edit script for horarios-141-2026-01-20.xlsx update
# Updates three sheets (LP1912, LP1912-215, 6203-6173) with the latest scrape
# snapshot (Última actualización: 09:31:25), re-sorted/new rows included.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value2 = "Última actualización: 09:31:25"

$ws.Range("A3").Value2 = "Total filas: 124"

$ws.Range("A67").Value2 = "08:20:43"
$ws.Range("B67").Value2 = "08:21"
$ws.Range("C67").Value2 = "215B_EL PATO"
$ws.Range("D67").Value2 = 1
$ws.Range("E67").Value2 = "LP1912"

$ws.Range("A68").Value2 = "06:44:15"
$ws.Range("B68").Value2 = "08:21"
$ws.Range("C68").Value2 = "26_HERNANDEZ"
$ws.Range("D68").Value2 = 97
$ws.Range("E68").Value2 = "LP1912"

$ws.Range("A99").Value2 = "09:31:25"
$ws.Range("B99").Value2 = "09:31"
$ws.Range("C99").Value2 = "23_HERNANDEZ"
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = "LP1912"

$ws.Range("A100").Value2 = "07:57:27"
$ws.Range("B100").Value2 = "09:32"
$ws.Range("C100").Value2 = "15_ABASTO"
$ws.Range("D100").Value2 = 95
$ws.Range("E100").Value2 = "LP1912"

$ws.Range("A101").Value2 = "07:57:27"
$ws.Range("B101").Value2 = "09:33"
$ws.Range("C101").Value2 = "10_OLMOS"
$ws.Range("D101").Value2 = 96
$ws.Range("E101").Value2 = "LP1912"

$ws.Range("A102").Value2 = "08:42:31"
$ws.Range("B102").Value2 = "09:34"
$ws.Range("C102").Value2 = "23_HERNANDEZ"
$ws.Range("D102").Value2 = 52
$ws.Range("E102").Value2 = "LP1912"

$ws.Range("A103").Value2 = "08:20:43"
$ws.Range("B103").Value2 = "09:41"
$ws.Range("C103").Value2 = "215C_EL PATO"
$ws.Range("D103").Value2 = 81
$ws.Range("E103").Value2 = "LP1912"

$ws.Range("A104").Value2 = "09:31:25"
$ws.Range("B104").Value2 = "09:42"
$ws.Range("C104").Value2 = "16_SANTA ANA"
$ws.Range("D104").Value2 = 11
$ws.Range("E104").Value2 = "LP1912"

$ws.Range("A105").Value2 = "07:57:27"
$ws.Range("B105").Value2 = "09:42"
$ws.Range("C105").Value2 = "215C_EL PATO"
$ws.Range("D105").Value2 = 105
$ws.Range("E105").Value2 = "LP1912"

$ws.Range("A106").Value2 = "07:57:27"
$ws.Range("B106").Value2 = "09:43"
$ws.Range("C106").Value2 = "14_ABASTO"
$ws.Range("D106").Value2 = 106
$ws.Range("E106").Value2 = "LP1912"

$ws.Range("A107").Value2 = "08:55:44"
$ws.Range("B107").Value2 = "09:52"
$ws.Range("C107").Value2 = "15_ABASTO"
$ws.Range("D107").Value2 = 57
$ws.Range("E107").Value2 = "LP1912"

$ws.Range("A108").Value2 = "08:55:44"
$ws.Range("B108").Value2 = "09:53"
$ws.Range("C108").Value2 = "10_OLMOS"
$ws.Range("D108").Value2 = 58
$ws.Range("E108").Value2 = "LP1912"

$ws.Range("A109").Value2 = "09:31:25"
$ws.Range("B109").Value2 = "09:59"
$ws.Range("C109").Value2 = "23_HERNANDEZ"
$ws.Range("D109").Value2 = 28
$ws.Range("E109").Value2 = "LP1912"

$ws.Range("A110").Value2 = "09:31:25"
$ws.Range("B110").Value2 = "10:02"
$ws.Range("C110").Value2 = "17_ROMERO"
$ws.Range("D110").Value2 = 31
$ws.Range("E110").Value2 = "LP1912"

$ws.Range("A111").Value2 = "09:31:25"
$ws.Range("B111").Value2 = "10:03"
$ws.Range("C111").Value2 = "11_ETCHEVERRY"
$ws.Range("D111").Value2 = 32
$ws.Range("E111").Value2 = "LP1912"

$ws.Range("A112").Value2 = "08:42:31"
$ws.Range("B112").Value2 = "10:07"
$ws.Range("C112").Value2 = "10_OLMOS"
$ws.Range("D112").Value2 = 85
$ws.Range("E112").Value2 = "LP1912"

$ws.Range("A113").Value2 = "08:20:43"
$ws.Range("B113").Value2 = "10:08"
$ws.Range("C113").Value2 = "10_OLMOS"
$ws.Range("D113").Value2 = 108
$ws.Range("E113").Value2 = "LP1912"

$ws.Range("A114").Value2 = "08:20:43"
$ws.Range("B114").Value2 = "10:12"
$ws.Range("C114").Value2 = "15_ABASTO"
$ws.Range("D114").Value2 = 112
$ws.Range("E114").Value2 = "LP1912"

$ws.Range("A115").Value2 = "09:31:25"
$ws.Range("B115").Value2 = "10:13"
$ws.Range("C115").Value2 = "10_OLMOS"
$ws.Range("D115").Value2 = 42
$ws.Range("E115").Value2 = "LP1912"

$ws.Range("A116").Value2 = "08:42:31"
$ws.Range("B116").Value2 = "10:21"
$ws.Range("C116").Value2 = "26_HERNANDEZ"
$ws.Range("D116").Value2 = 99
$ws.Range("E116").Value2 = "LP1912"

$ws.Range("A117").Value2 = "09:31:25"
$ws.Range("B117").Value2 = "10:22"
$ws.Range("C117").Value2 = "16_SANTA ANA"
$ws.Range("D117").Value2 = 51
$ws.Range("E117").Value2 = "LP1912"

$ws.Range("A118").Value2 = "09:31:25"
$ws.Range("B118").Value2 = "10:23"
$ws.Range("C118").Value2 = "11_ETCHEVERRY"
$ws.Range("D118").Value2 = 52
$ws.Range("E118").Value2 = "LP1912"

$ws.Range("A119").Value2 = "08:42:31"
$ws.Range("B119").Value2 = "10:26"
$ws.Range("C119").Value2 = "215A_EL PATO"
$ws.Range("D119").Value2 = 104
$ws.Range("E119").Value2 = "LP1912"

$ws.Range("A120").Value2 = "08:55:44"
$ws.Range("B120").Value2 = "10:41"
$ws.Range("C120").Value2 = "17_ROMERO"
$ws.Range("D120").Value2 = 106
$ws.Range("E120").Value2 = "LP1912"

$ws.Range("A121").Value2 = "09:31:25"
$ws.Range("B121").Value2 = "10:42"
$ws.Range("C121").Value2 = "17_ROMERO"
$ws.Range("D121").Value2 = 71
$ws.Range("E121").Value2 = "LP1912"

$ws.Range("A122").Value2 = "08:55:44"
$ws.Range("B122").Value2 = "10:43"
$ws.Range("C122").Value2 = "14_ABASTO"
$ws.Range("D122").Value2 = 108
$ws.Range("E122").Value2 = "LP1912"

$ws.Range("A123").Value2 = "09:31:25"
$ws.Range("B123").Value2 = "10:59"
$ws.Range("C123").Value2 = "27_EL RETIRO"
$ws.Range("D123").Value2 = 88
$ws.Range("E123").Value2 = "LP1912"

$ws.Range("A124").Value2 = "09:31:25"
$ws.Range("B124").Value2 = "11:02"
$ws.Range("C124").Value2 = "215C_EL PATO"
$ws.Range("D124").Value2 = 91
$ws.Range("E124").Value2 = "LP1912"

$ws.Range("A125").Value2 = "09:31:25"
$ws.Range("B125").Value2 = "11:17"
$ws.Range("C125").Value2 = "16_P MOR-167 Y 521"
$ws.Range("D125").Value2 = 106
$ws.Range("E125").Value2 = "LP1912"

$ws.Range("A126").Value2 = "09:31:25"
$ws.Range("B126").Value2 = "11:19"
$ws.Range("C126").Value2 = "86_EST CHICA-ESC AGRARIA"
$ws.Range("D126").Value2 = 108
$ws.Range("E126").Value2 = "LP1912"

$ws.Range("A127").Value2 = "09:31:25"
$ws.Range("B127").Value2 = "11:21"
$ws.Range("C127").Value2 = "26_HERNANDEZ"
$ws.Range("D127").Value2 = 110
$ws.Range("E127").Value2 = "LP1912"

$ws.Range("A128").Value2 = "09:31:25"
$ws.Range("B128").Value2 = "11:26"
$ws.Range("C128").Value2 = "16_P MOR-SANTA ANA"
$ws.Range("D128").Value2 = 115
$ws.Range("E128").Value2 = "LP1912"

$ws.Range("A129").Value2 = "09:31:25"
$ws.Range("B129").Value2 = "11:27"
$ws.Range("C129").Value2 = "225_C ROCA-H SUR"
$ws.Range("D129").Value2 = 116
$ws.Range("E129").Value2 = "LP1912"


$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Range("A2").Value2 = "Última actualización: 09:31:25"

$ws.Range("A3").Value2 = "Total filas: 15"

$ws.Range("A20").Value2 = "09:31:25"
$ws.Range("B20").Value2 = "11:02"
$ws.Range("C20").Value2 = "215C_EL PATO"
$ws.Range("D20").Value2 = 91
$ws.Range("E20").Value2 = "LP1912"


$ws = $wb.Worksheets.Item("6203-6173")

$ws.Range("A2").Value2 = "Última actualización: 09:31:25"

$ws.Range("A3").Value2 = "Total filas: 24"

$ws.Range("A29").Value2 = "09:31:25"
$ws.Range("B29").Value2 = "11:14"
$ws.Range("C29").Value2 = "215C_LA PLATA"
$ws.Range("D29").Value2 = 103
$ws.Range("E29").Value2 = "L6203"

